$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Execute" (column C) values for all rows except the
# "Extension Payments" block (rows 21-27), matching the RAD data
# refresh for the EmailNoMatch test case.
$ws.Range("C2:C20").Clear()
$ws.Range("C28:C54").Clear()

# Refresh the execution timestamps for the Extension Payments rows
# (21-27) with the new run's values.
$ws.Range("B21").Value = "Wed Mar 20 23:03:47 EDT 2024"
$ws.Range("B22").Value = "Wed Mar 20 23:03:59 EDT 2024"
$ws.Range("B23").Value = "Wed Mar 20 23:04:11 EDT 2024"
$ws.Range("B24").Value = "Wed Mar 20 23:04:23 EDT 2024"
$ws.Range("B25").Value = "Wed Mar 20 23:04:35 EDT 2024"
$ws.Range("B26").Value = "Wed Mar 20 23:04:47 EDT 2024"
$ws.Range("B27").Value = "Wed Mar 20 23:04:59 EDT 2024"

# Update the active selection to reflect where the user left off.
$ws.Range("C48").Select()
